$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <= original row 5 data
$ws.Range("D2").Value = 44904
$ws.Range("K2").Value = 'Castle Brite'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1550
$ws.Range("T2").Value = 10

# Row 3 <= original row 6 data
$ws.Range("D3").Value = 44904
$ws.Range("K3").Value = 'Castle Brite'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 14000
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1400
$ws.Range("T3").Value = 10

# Row 4 <= original row 12 data
$ws.Range("D4").Value = 44537
$ws.Range("K4").Value = 'Castle Brite'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 21500
$ws.Range("P4").Value = 21250
$ws.Range("Q4").Value = '$/caja 15 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1417
$ws.Range("T4").Value = 15

# Row 5 <= original row 13 data
$ws.Range("D5").Value = 44189
$ws.Range("K5").Value = 'Dina'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 16562
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 920
$ws.Range("T5").Value = 18

# Row 6 <= original row 3 data
$ws.Range("D6").Value = 44552
$ws.Range("K6").Value = 'Castle Brite'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 15500
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15750
$ws.Range("Q6").Value = '$/caja 15 kilos'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 1050
$ws.Range("T6").Value = 15

# Row 7 <= original row 8 data
$ws.Range("D7").Value = 44907
$ws.Range("K7").Value = 'Castle Brite'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 15500
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1550
$ws.Range("T7").Value = 10

# Row 8 <= original row 9 data
$ws.Range("D8").Value = 44907
$ws.Range("K8").Value = 'Castle Brite'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1400
$ws.Range("T8").Value = 10

# Row 9 <= original row 14 data
$ws.Range("D9").Value = 44176
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 17400
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 967
$ws.Range("T9").Value = 18

# Row 10 <= original row 16 data
$ws.Range("D10").Value = 44551
$ws.Range("K10").Value = 'Castle Brite'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 15500
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15750
$ws.Range("Q10").Value = '$/caja 15 kilos'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1050
$ws.Range("T10").Value = 15

# Row 12 <= original row 7 data
$ws.Range("D12").Value = 44187
$ws.Range("K12").Value = 'Dina'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 55
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15455
$ws.Range("Q12").Value = '$/caja 15 kilos granel'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 1030
$ws.Range("T12").Value = 15

# Row 13 <= original row 10 data
$ws.Range("D13").Value = 44181
$ws.Range("K13").Value = 'Modesto'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Región de Coquimbo'
$ws.Range("S13").Value = 1139
$ws.Range("T13").Value = 18

# Row 14 <= original row 4 data
$ws.Range("D14").Value = 44165
$ws.Range("K14").Value = 'Castle Brite'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 17000
$ws.Range("P14").Value = 16500
$ws.Range("Q14").Value = '$/caja 15 kilos granel'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 1100
$ws.Range("T14").Value = 15

# Row 15 <= original row 2 data
$ws.Range("D15").Value = 44174
$ws.Range("K15").Value = 'Castle Brite'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = 9000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 9467
$ws.Range("Q15").Value = '$/caja 10 kilos'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 947
$ws.Range("T15").Value = 10

# Row 16 <= original row 15 data
$ws.Range("D16").Value = 44168
$ws.Range("K16").Value = 'Castle Brite'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 16000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 16500
$ws.Range("Q16").Value = '$/caja 16 kilos granel'
$ws.Range("R16").Value = 'Región de Coquimbo'
$ws.Range("S16").Value = 1031
$ws.Range("T16").Value = 16
